# The commit swaps the two embedded theme color schemes:
#   ppt/theme/theme1.xml (the deck's active "Simple Light" theme) picks up
#   the color values that used to live in ppt/theme/theme2.xml (the
#   "Default" theme used by the notes master) -- font scheme and format
#   scheme are identical between the two themes, only the 12 clrScheme
#   colors actually change.
#
# Re-create that by pushing the "Default" theme's RGB values onto the
# active presentation theme's ThemeColorScheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink -- indices 1..12 in that order). PowerPoint's
# RGB property packs colors as 0xBBGGRR (the classic VBA RGB() order), so
# a target hex RRGGBB of e.g. 158158 is written as 0x588115.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x588115   # dk2      -> 158158
$tcs.Item(4).RGB  = 0xF3F3F3   # lt2      -> F3F3F3
$tcs.Item(5).RGB  = 0xC78D05   # accent1  -> 058DC7
$tcs.Item(6).RGB  = 0x32B450   # accent2  -> 50B432
$tcs.Item(7).RGB  = 0x1B56ED   # accent3  -> ED561B
$tcs.Item(8).RGB  = 0x00EFED   # accent4  -> EDEF00
$tcs.Item(9).RGB  = 0xE5CB24   # accent5  -> 24CBE5
$tcs.Item(10).RGB = 0x72E564   # accent6  -> 64E572
$tcs.Item(11).RGB = 0xCC0022   # hlink    -> 2200CC
$tcs.Item(12).RGB = 0x8B1A55   # folHlink -> 551A8B
